$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 with the new master-server record.
# Insert shared strings in the same order the target workbook uses
# (ServerID, then IP, then Name/ID) so the sharedStrings table indices line up.
$ws.Range("B2").Value = "000106001"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "MasterServer_1"

# C2 has no pre-existing style (unlike A2/B2, which already use the "Text" format),
# so explicitly give it the same text number format before writing its value.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "MasterServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 2001

# Move the selection to H3, matching the saved view state.
$ws.Range("H3").Select() | Out-Null
